$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text labels -----------------------------------------------------------
# Written in the order that reproduces the original sharedStrings table
# order (Governo, receita, gasto, imposto_renda, taxa_juros, deposito_compulsório)
$ws.Range("A1").Value = "Governo"
$ws.Range("A3").Value = "receita"
$ws.Range("A2").Value = "gasto"
$ws.Range("A4").Value = "imposto_renda"
$ws.Range("A5").Value = "taxa_juros"
$ws.Range("A6").Value = "deposito_compulsório"

# --- Numeric values ----------------------------------------------------------
$ws.Range("B2").Value = 100000
$ws.Range("B3").Value = 150000
$ws.Range("B4").Value = 0.1
$ws.Range("B5").Value = 0.04
$ws.Range("B6").Value = 0.5

# --- Number formats ----------------------------------------------------------
# Currency ("Moeda") format for gasto/receita
$ws.Range("B2:B3").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"
# Percentage format for the rate rows
$ws.Range("B4:B6").NumberFormat = "0%"

# --- Header row formatting: bold, centered, merged A1:B1 --------------------
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").HorizontalAlignment = -4108
$ws.Range("A1:B1").Merge()

# --- Column widths -------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.44140625
$ws.Columns.Item(2).ColumnWidth = 12.44140625

# --- Sheet view: zoom + selection -------------------------------------------
$excel.ActiveWindow.Zoom = 160
[void]$ws.Range("C8").Select()
